$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns I and J (teacher collection / price columns) mirror the
# existing column width (14.5) used by columns A-H.
$ws.Columns("I:J").ColumnWidth = 13.66

# Give every new I/J cell the same bordered, unfilled, default-font look
# already used elsewhere in the sheet (reuses the existing thin grey
# border definition).
$ws.Range("I1:J53").Borders.ColorIndex = 3

# J10 becomes a dynamically calculated (time-formatted) price/value cell,
# matching the sheet's banded white fill.
$ws.Range("J10").NumberFormat = "h:mm AM/PM"
$ws.Range("J10").Interior.ColorIndex = 2
